# Update the probability-threshold table on Sheet1 to the new Zn/2His_1Glu
# values and leave the live selection on C3 (matching the author's last
# on-screen selection when the file was uploaded).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# alpha_distance_range (row 2) and beta_distance_range (row 3): Min/Max
$ws.Range("B2").Value = 5.7
$ws.Range("C2").Value = 10.7
$ws.Range("B3").Value = 5.7
$ws.Range("C3").Value = 9.3000000000000007

# Restore the active cell/selection to C3, as captured in the saved view.
$ws.Activate() | Out-Null
$ws.Range("C3").Select() | Out-Null
